$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 41 (G41=5478)
$ws.Range("H41").Value = 1487.4706
$ws.Range("I41").Value = 1832.5555
$ws.Range("K41").Value = 1832.5555
$ws.Range("M41").Value = -1392.5555

# Row 51 (G51=5486)
$ws.Range("H51").Value = 11117.333
$ws.Range("I51").Value = 18550.166
$ws.Range("J51").Value = 3684.5
$ws.Range("K51").Value = 18550.166
$ws.Range("L51").Value = 3684.5
$ws.Range("M51").Value = -18066.166
$ws.Range("N51").Value = -4652.5

# Row 92 (G92=19901)
$ws.Range("H92").Value = 581.3684
$ws.Range("I92").Value = 665.7857
$ws.Range("J92").Value = 345
$ws.Range("K92").Value = 665.7857
$ws.Range("L92").Value = 345
$ws.Range("M92").Value = 582.2143
$ws.Range("N92").Value = -2841

# Row 106 (G106=19903)
$ws.Range("H106").Value = 2195.2
$ws.Range("I106").Value = 2195.2
$ws.Range("K106").Value = 2195.2
$ws.Range("M106").Value = -1564.2

# Row 129 (G129=36115)
$ws.Range("H129").Value = 2896.54
$ws.Range("I129").Value = 11608.667
$ws.Range("K129").Value = 34826.001
$ws.Range("M129").Value = -29826.001

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (G32=44147)
$ws.Range("H32").Value = 26810.646
$ws.Range("I32").Value = 4449.0527
$ws.Range("K32").Value = 4449.0527
$ws.Range("M32").Value = -4162.0527

# Row 110 (G110=27708)
$ws.Range("H110").Value = 23858922
$ws.Range("I110").Value = 38539616
$ws.Range("J110").Value = 2789
$ws.Range("K110").Value = 38539616
$ws.Range("L110").Value = 2789
$ws.Range("M110").Value = -38537571
$ws.Range("N110").Value = -6879

# Row 122 (G122=36168)
$ws.Range("H122").Value = 1869.25
$ws.Range("I122").Value = 1761.409
$ws.Range("J122").Value = 2264.6667
$ws.Range("K122").Value = 5284.227000000001
$ws.Range("L122").Value = 6794.000100000001
$ws.Range("M122").Value = -2834.227000000001
$ws.Range("N122").Value = -11694.0001

# Row 132 (G132=43997)
$ws.Range("H132").Value = 6405.2925
$ws.Range("I132").Value = 6429.057
$ws.Range("J132").Value = 6266.6665
$ws.Range("K132").Value = 19287.171
$ws.Range("L132").Value = 18799.9995
$ws.Range("M132").Value = -16757.171
$ws.Range("N132").Value = -23859.9995

$ws = $wb.Worksheets.Item("BSM")
# Row 41 (G41=22899)
$ws.Range("H41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()

# Row 48 (G48=22888)
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()

# Row 86 (G86=12526)
$ws.Range("H86").Value = 94407.164
$ws.Range("I86").Value = 139699.88
$ws.Range("J86").Value = 3821.75
$ws.Range("K86").Value = 139699.88
$ws.Range("L86").Value = 3821.75
$ws.Range("M86").Value = -138576.88
$ws.Range("N86").Value = -6067.75

# Row 89 (G89=12526)
$ws.Range("H89").Value = 94407.164
$ws.Range("I89").Value = 139699.88
$ws.Range("J89").Value = 3821.75
$ws.Range("K89").Value = 698499.4
$ws.Range("L89").Value = 19108.75
$ws.Range("M89").Value = -692883.4
$ws.Range("N89").Value = -30340.75

# Row 94 (G94=19939)
$ws.Range("H94").Value = 111403.336
$ws.Range("I94").Value = 125266.25
$ws.Range("K94").Value = 125266.25
$ws.Range("M94").Value = -124815.25

# Row 105 (G105=19947)
$ws.Range("H105").Value = 64280.656
$ws.Range("I105").Value = 41639.6
$ws.Range("J105").Value = 145141.58
$ws.Range("K105").Value = 41639.6
$ws.Range("L105").Value = 145141.58
$ws.Range("M105").Value = -39892.6
$ws.Range("N105").Value = -148635.58

$ws = $wb.Worksheets.Item("CRP")
# Row 16 (G16=27691)
$ws.Range("H16").Value = 977.625
$ws.Range("I16").Value = 505.5
$ws.Range("J16").Value = 1135
$ws.Range("K16").Value = 505.5
$ws.Range("L16").Value = 1135
$ws.Range("M16").Value = -218.5
$ws.Range("N16").Value = -1709

# Row 113 (G113=27691)
$ws.Range("H113").Value = 977.625
$ws.Range("I113").Value = 505.5
$ws.Range("J113").Value = 1135
$ws.Range("K113").Value = 505.5
$ws.Range("L113").Value = 1135
$ws.Range("M113").Value = 1664.5
$ws.Range("N113").Value = -5475

# Row 122 (G122=36196)
$ws.Range("H122").Value = 2188.4412
$ws.Range("I122").Value = 2099.75
$ws.Range("J122").Value = 2602.3333
$ws.Range("K122").Value = 6299.25
$ws.Range("L122").Value = 7806.999899999999
$ws.Range("M122").Value = -3849.25
$ws.Range("N122").Value = -12706.9999

$ws = $wb.Worksheets.Item("CUL")
# Row 68 (G68=12895)
$ws.Range("H68").Value = 1942.8154
$ws.Range("I68").Value = 1294.1538
$ws.Range("J68").Value = 2375.2563
$ws.Range("K68").Value = 3882.4614
$ws.Range("L68").Value = 7125.7689
$ws.Range("M68").Value = -3071.4614
$ws.Range("N68").Value = -8747.768899999999

# Row 71 (G71=12895)
$ws.Range("H71").Value = 1942.8154
$ws.Range("I71").Value = 1294.1538
$ws.Range("J71").Value = 2375.2563
$ws.Range("K71").Value = 11647.3842
$ws.Range("L71").Value = 21377.3067
$ws.Range("M71").Value = -7591.3842
$ws.Range("N71").Value = -29489.3067

# Row 122 (G122=36078)
$ws.Range("H122").Value = 6574.706
$ws.Range("I122").Value = 253.14285
$ws.Range("J122").Value = 10999.8
$ws.Range("K122").Value = 2278.28565
$ws.Range("L122").Value = 98998.2
$ws.Range("M122").Value = 171.7143499999997
$ws.Range("N122").Value = -103898.2

# Row 131 (G131=36060)
$ws.Range("H131").Value = 1534.125
$ws.Range("J131").Value = 1528.4698
$ws.Range("L131").Value = 4585.4094
$ws.Range("N131").Value = -14665.4094

# Row 137 (G137=44088)
$ws.Range("H137").Value = 11942520
$ws.Range("I137").Value = 68166.266
$ws.Range("J137").Value = 25643698
$ws.Range("K137").Value = 204498.798
$ws.Range("L137").Value = 76931094
$ws.Range("M137").Value = -199398.798
$ws.Range("N137").Value = -76941294

$ws = $wb.Worksheets.Item("GSM")
# Row 97 (G97=19940)
$ws.Range("H97").Value = 45456590
$ws.Range("I97").Value = 58825664
$ws.Range("K97").Value = 58825664
$ws.Range("M97").Value = -58825168

# Row 122 (G122=36182)
$ws.Range("H122").Value = 3224.2273
$ws.Range("I122").Value = 3554.6924
$ws.Range("J122").Value = 2746.889
$ws.Range("K122").Value = 10664.0772
$ws.Range("L122").Value = 8240.667000000001
$ws.Range("M122").Value = -8214.0772
$ws.Range("N122").Value = -13140.667

# Row 134 (G134=42064)
$ws.Range("H134").Value = 22316.8
$ws.Range("J134").Value = 22316.8
$ws.Range("L134").Value = 66950.39999999999
$ws.Range("N134").Value = -72020.39999999999

$ws = $wb.Worksheets.Item("LTW")
# Row 7 (G7=36249)
$ws.Range("H7").Value = 1694.5714
$ws.Range("I7").Value = 1241.4375
$ws.Range("J7").Value = 2298.75
$ws.Range("K7").Value = 1241.4375
$ws.Range("L7").Value = 2298.75
$ws.Range("M7").Value = -1129.4375
$ws.Range("N7").Value = -2522.75

# Row 16 (G16=5289)
$ws.Range("H16").Value = 7876985.5
$ws.Range("I16").Value = 14000943
$ws.Range("J16").Value = 3326
$ws.Range("K16").Value = 14000943
$ws.Range("L16").Value = 3326
$ws.Range("M16").Value = -14000773
$ws.Range("N16").Value = -3666

# Row 61 (G61=27740)
$ws.Range("H61").Value = 1864.7778
$ws.Range("I61").Value = 1525.375
$ws.Range("J61").Value = 4580
$ws.Range("K61").Value = 1525.375
$ws.Range("L61").Value = 4580
$ws.Range("M61").Value = -1323.375
$ws.Range("N61").Value = -4984

# Row 113 (G113=27740)
$ws.Range("H113").Value = 1864.7778
$ws.Range("I113").Value = 1525.375
$ws.Range("J113").Value = 4580
$ws.Range("K113").Value = 1525.375
$ws.Range("L113").Value = 4580
$ws.Range("M113").Value = 644.625
$ws.Range("N113").Value = -8920

# Row 122 (G122=36247)
$ws.Range("H122").Value = 3448.875
$ws.Range("I122").Value = 3370.1428
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 10110.4284
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -7660.428400000001
$ws.Range("N122").Value = -16900

# Row 126 (G126=36249)
$ws.Range("H126").Value = 1694.5714
$ws.Range("I126").Value = 1241.4375
$ws.Range("J126").Value = 2298.75
$ws.Range("K126").Value = 3724.3125
$ws.Range("L126").Value = 6896.25
$ws.Range("M126").Value = -1254.3125
$ws.Range("N126").Value = -11836.25

# Row 134 (G134=42024)
$ws.Range("H134").Value = 60757.668
$ws.Range("J134").Value = 60757.668
$ws.Range("L134").Value = 60757.668
$ws.Range("N134").Value = -70897.66800000001
